# Add new columns I (I0) and J (IF) to the worksheet, as per commit "I0 and IF added"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - set values then copy formatting (style) from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-19: columns I and J hold plain numbers (no special style)
$data = @{
  2  = @(7, 8)
  3  = @(9, 9)
  4  = @(7, 7)
  5  = @(4, 5)
  6  = @(7, 8)
  7  = @(5, 5)
  8  = @(5, 6)
  9  = @(5, 5)
  10 = @(8, 8)
  11 = @(6, 8)
  12 = @(6, 6)
  13 = @(8, 8)
  14 = @(9, 9)
  15 = @(9, 9)
  16 = @(6, 6)
  17 = @(8, 8)
  18 = @(7, 7)
  19 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value  = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
